$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the numeric-looking Price cells to stay text (matches source inlineStr cells)
$textCells = @("D2", "D3", "D5", "D6", "D8", "D9", "D11", "D12", "D13", "D14", "D16", "D17", "D18", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D27", "D29", "D30", "D31", "D32", "D33", "D34", "D36", "D38", "D39", "D40", "D42", "D43", "D44", "D46", "D47", "D49", "D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '59.827.22'
$ws.Range('E2').Value = '  -5.66%  '
$ws.Range('D3').Value = '2.972.38'
$ws.Range('E3').Value = '  -6.31%  '
$ws.Range('E4').Value = '  +0.21%  '
$ws.Range('D5').Value = '565.60'
$ws.Range('E5').Value = '  -6.06%  '
$ws.Range('D6').Value = '123.67'
$ws.Range('E6').Value = '  -9.12%  '
$ws.Range('E7').Value = '  +0.21%  '
$ws.Range('D8').Value = '2.964.31'
$ws.Range('E8').Value = '  -6.52%  '
$ws.Range('D9').Value = '0.500'
$ws.Range('E9').Value = '  -2.49%  '
$ws.Range('E10').Value = '  -9.37%  '
$ws.Range('D11').Value = '4.92'
$ws.Range('E11').Value = '  -8.39%  '
$ws.Range('D12').Value = '0.439'
$ws.Range('E12').Value = '  -3.55%  '
$ws.Range('D13').Value = '0.0000218'
$ws.Range('E13').Value = '  -9.11%  '
$ws.Range('D14').Value = '32.15'
$ws.Range('E14').Value = '  -7.77%  '
$ws.Range('E15').Value = '  -0.13%  '
$ws.Range('D16').Value = '3.470.69'
$ws.Range('E16').Value = '  -6.09%  '
$ws.Range('D17').Value = '2.972.72'
$ws.Range('E17').Value = '  -6.27%  '
$ws.Range('D18').Value = '59.883.92'
$ws.Range('E18').Value = '  -5.54%  '
$ws.Range('D19').Value = '6.48'
$ws.Range('E19').Value = '  -1.65%  '
$ws.Range('D20').Value = '425.26'
$ws.Range('E20').Value = '  -8.05%  '
$ws.Range('D21').Value = '13.08'
$ws.Range('E21').Value = '  -6.49%  '
$ws.Range('D22').Value = '0.665'
$ws.Range('E22').Value = '  -4.73%  '
$ws.Range('D23').Value = '6.95'
$ws.Range('E23').Value = '  -9.53%  '
$ws.Range('D24').Value = '12.76'
$ws.Range('E24').Value = '  -4.10%  '
$ws.Range('D25').Value = '78.96'
$ws.Range('E25').Value = '  -5.11%  '
$ws.Range('E26').Value = '  +0.08%  '
$ws.Range('D27').Value = '1.00'
$ws.Range('E27').Value = '  +0.08%  '
$ws.Range('E28').Value = '  -7.35%  '
$ws.Range('D29').Value = '1.92'
$ws.Range('E29').Value = '  -7.69%  '
$ws.Range('D30').Value = '7.12'
$ws.Range('E30').Value = '  -7.67%  '
$ws.Range('D31').Value = '6.04'
$ws.Range('E31').Value = '  -11.10%  '
$ws.Range('D32').Value = '24.97'
$ws.Range('E32').Value = '  -8.05%  '
$ws.Range('D33').Value = '0.0959'
$ws.Range('E33').Value = '  -4.87%  '
$ws.Range('D34').Value = '5.53'
$ws.Range('E34').Value = '  -6.29%  '
$ws.Range('E35').Value = '  -2.15%  '
$ws.Range('D36').Value = '0.916'
$ws.Range('E36').Value = '  -10.17%  '
$ws.Range('E37').Value = '  -19.26%  '
$ws.Range('D38').Value = '8.43'
$ws.Range('E38').Value = '  +3.60%  '
$ws.Range('D39').Value = '0.0₃0641'
$ws.Range('E39').Value = '  -12.61%  '
$ws.Range('D40').Value = '0.0352'
$ws.Range('E40').Value = '  -9.81%  '
$ws.Range('E41').Value = '  -6.12%  '
$ws.Range('D42').Value = '2.654.08'
$ws.Range('E42').Value = '  -5.51%  '
$ws.Range('D43').Value = '362.52'
$ws.Range('E43').Value = '  -8.04%  '
$ws.Range('D44').Value = '2.39'
$ws.Range('E44').Value = '  -9.27%  '
$ws.Range('D46').Value = '120.15'
$ws.Range('E46').Value = '  -4.82%  '
$ws.Range('D47').Value = '0.231'
$ws.Range('E47').Value = '  -8.11%  '
$ws.Range('E48').Value = '  -3.92%  '
$ws.Range('D49').Value = '1.95'
$ws.Range('E49').Value = '  -7.94%  '
$ws.Range('D50').Value = '23.07'
$ws.Range('E50').Value = '  -8.42%  '
$ws.Range('E51').Value = '  -8.51%  '

# Reset formatting/style back to the sheet default so no stray styles are introduced
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "General"
    $ws.Range($addr).Style = "Normal"
}
